# Apply changes described by the diff:
# 1. Insert a new "Player Info" sheet as the first sheet with player data.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling" sheets,
#    and change the URL data value to just the match code number.

$wb = $excel.ActiveWorkbook

# --- Update "ODI Batting" sheet (currently sheet1, header D1) ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
# Leading apostrophe forces the numeric-looking match code to stay text
# (matches the data convention used throughout this workbook), then
# ClearFormats drops the transient "number stored as text" style marker.
$batting.Range("D2").Value = "'4267"
$batting.Range("D2").ClearFormats()

# --- Update "ODI Bowling" sheet (currently sheet2, header B1) ---
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4267"
$bowling.Range("B2").ClearFormats()

# --- Insert new "Player Info" sheet before "ODI Batting" (becomes first sheet) ---
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Re-fetch the "ODI Batting" sheet reference (it may have shifted after the
# insert above) before copying its header formatting to the new sheet.
$battingRef = $wb.Worksheets.Item("ODI Batting")
$battingRef.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A2").Value = "'4782"
$playerInfo.Range("A2").ClearFormats()
$playerInfo.Range("B2").Value = "Zahir Khan Pakteen"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Wrist Spin (Chinaman)"
